$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure cells are treated as plain text so values such as "205.98" or
# "1.564.20" are not auto-converted into numbers/dates by Excel.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '26.868.12'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  -1.08%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.564.20'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +0.27%  '
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  -0.18%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '205.98'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  -0.21%  '
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -0.99%  '
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  -0.16%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '21.82'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -1.69%  '
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -0.19%  '
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -1.24%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0864'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +0.26%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.785.59'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +0.19%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.574.64'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +0.96%  '
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  -1.26%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.514'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -0.23%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '26.871.94'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  -1.06%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '61.32'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -2.35%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '215.22'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +0.55%  '
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +1.79%  '
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  -0.56%  '
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -0.15%  '
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +0.36%  '
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -1.70%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.01'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +1.53%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '153.46'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +0.87%  '
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +1.55%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '14.94'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  +0.34%  '
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -0.19%  '
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  -1.00%  '
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +0.87%  '
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -3.67%  '
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -0.10%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.401.97'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +1.28%  '
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +0.02%  '
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  -0.84%  '
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -0.35%  '
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  -3.96%  '
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -0.73%  '
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  +2.15%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.813'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  +0.30%  '
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  -0.14%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.995'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +1.35%  '
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +6.57%  '
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +1.10%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '63.60'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  +0.49%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.699.36'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +0.29%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '86.58'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +1.33%  '
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +1.85%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0₇0973'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -1.18%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0951'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +0.99%  '
